$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 290 (shifts existing rows 290-309 down to 293-312)
$ws.Range("A290:T292").EntireRow.Insert()

# Fill in the 3 new rows (290-292) with new Murcott / Region de O'Higgins data
# Row 290
$ws.Range("A290").Value2 = 4
$ws.Range("B290").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C290").Value2 = "Los Lagos"
$ws.Range("D290").Value2 = 44826
$ws.Range("E290").Value2 = 10
$ws.Range("F290").Value2 = "Fruta"
$ws.Range("G290").Value2 = 100102
$ws.Range("H290").Value2 = "Cítricos"
$ws.Range("I290").Value2 = 100102004
$ws.Range("J290").Value2 = "Mandarina"
$ws.Range("K290").Value2 = "Murcott"
$ws.Range("L290").Value2 = "Especial"
$ws.Range("M290").Value2 = 150
$ws.Range("N290").Value2 = 11000
$ws.Range("O290").Value2 = 11000
$ws.Range("P290").Value2 = 11000
$ws.Range("Q290").Value2 = "$/bandeja 10 kilos"
$ws.Range("R290").Value2 = "Región de O'Higgins"
$ws.Range("S290").Value2 = 1100
$ws.Range("T290").Value2 = 10

# Row 291
$ws.Range("A291").Value2 = 4
$ws.Range("B291").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C291").Value2 = "Los Lagos"
$ws.Range("D291").Value2 = 44826
$ws.Range("E291").Value2 = 10
$ws.Range("F291").Value2 = "Fruta"
$ws.Range("G291").Value2 = 100102
$ws.Range("H291").Value2 = "Cítricos"
$ws.Range("I291").Value2 = 100102004
$ws.Range("J291").Value2 = "Mandarina"
$ws.Range("K291").Value2 = "Murcott"
$ws.Range("L291").Value2 = "Primera"
$ws.Range("M291").Value2 = 150
$ws.Range("N291").Value2 = 9500
$ws.Range("O291").Value2 = 9500
$ws.Range("P291").Value2 = 9500
$ws.Range("Q291").Value2 = "$/bandeja 10 kilos"
$ws.Range("R291").Value2 = "Región de O'Higgins"
$ws.Range("S291").Value2 = 950
$ws.Range("T291").Value2 = 10

# Row 292
$ws.Range("A292").Value2 = 4
$ws.Range("B292").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C292").Value2 = "Los Lagos"
$ws.Range("D292").Value2 = 44826
$ws.Range("E292").Value2 = 10
$ws.Range("F292").Value2 = "Fruta"
$ws.Range("G292").Value2 = 100102
$ws.Range("H292").Value2 = "Cítricos"
$ws.Range("I292").Value2 = 100102004
$ws.Range("J292").Value2 = "Mandarina"
$ws.Range("K292").Value2 = "Murcott"
$ws.Range("L292").Value2 = "Segunda"
$ws.Range("M292").Value2 = 150
$ws.Range("N292").Value2 = 7500
$ws.Range("O292").Value2 = 7500
$ws.Range("P292").Value2 = 7500
$ws.Range("Q292").Value2 = "$/bandeja 10 kilos"
$ws.Range("R292").Value2 = "Región de O'Higgins"
$ws.Range("S292").Value2 = 750
$ws.Range("T292").Value2 = 10
